# Populate the OctoberRaw sheet with this month's raw circulation data
# (values pulled in from the source system for October), then let the
# October sheet (which reads via =OctoberRaw!..) and the Yearly total
# sheet (which sums all 12 months) recalculate automatically.
$wb = $excel.ActiveWorkbook
$octRaw = $wb.Worksheets.Item("OctoberRaw")

# Header row
$octRaw.Cells.Item(1,1).Value = "Library"
$octRaw.Cells.Item(1,2).Value = "Items owned by this library checked out at this library this month"
$octRaw.Cells.Item(1,3).Value = "Items owned by other libraries checked out at this library this month"
$octRaw.Cells.Item(1,4).Value = "Total circulation this month"

# Column A library names, rows 2-54
$libraryNames = @(
    'Atchison Public Library'
    'Baldwin City Public Library'
    'Basehor Community Library'
    'Bern Community Library'
    'Bonner Springs City Library'
    'Burlingame Community Library'
    'Carbondale City Library'
    'Centralia Community Library'
    'Corning City Library'
    'Digital Content'
    'Doniphan County Library - Elwood'
    'Doniphan County Library - Highland'
    'Doniphan County Library - Troy'
    'Doniphan County Library - Wathena'
    'Effingham Community Library'
    'Eudora Community Library'
    'Everest, Barnes Reading Room'
    'Hiawatha, Morrill Public Library'
    'Highland Community College'
    'Holton, Beck-Bookman Library'
    'Horton Public Library'
    'Lansing Community Library'
    'Leavenworth Public Library'
    'Linwood Community Library'
    'Louisburg Library'
    'Lyndon Carnegie Library'
    'McLouth Public Library'
    'Meriden-Ozawkie Public Library'
    'Northeast Kansas Library System'
    'Nortonville Public Library'
    'Osage City Library'
    'Osawatomie Public Library'
    'Oskaloosa Public Library'
    'Ottawa Library'
    'Overbrook Public Library'
    'Paola Free Library'
    'Perry-Lecompton Community Library'
    'Pomona Community Library'
    'Prairie Hills Schools - Axtell Public School'
    'Prairie Hills Schools - Sabetha Elementary School'
    'Prairie Hills Schools - Sabetha High School'
    'Prairie Hills Schools - Sabetha Middle School'
    'Richmond Public Library'
    'Rossville Community Library'
    'Sabetha, Mary Cotton Library'
    'Seneca Free Library'
    'Silver Lake Library'
    'Tonganoxie Public Library'
    'Valley Falls, Delaware Township Library'
    'Wellsville City Library'
    'Wetmore Public Library'
    'Williamsburg Community Library'
    'Winchester Public Library'
)
for ($i = 0; $i -lt $libraryNames.Length; $i++) {
    $octRaw.Cells.Item($i + 2, 1).Value = $libraryNames[$i]
}

# Raw monthly totals (B=owned-here, C=owned-elsewhere, D=total), rows 2-54.
# Rows 11 (Digital Content) and 26 (Louisburg Library) are section/
# subtotal-free rows with no figures, so they are left blank.
$rawData = @{
    2 = @(4143, 1655, 5798)
    3 = @(2526, 570, 3096)
    4 = @(7249, 1051, 8300)
    5 = @(77, 31, 108)
    6 = @(5040, 1320, 6360)
    7 = @(497, 214, 711)
    8 = @(445, 182, 627)
    9 = @(287, 84, 371)
    10 = @(35, 1, 36)
    12 = @(77, 26, 103)
    13 = @(132, 68, 200)
    14 = @(413, 219, 632)
    15 = @(364, 85, 449)
    16 = @(166, 56, 222)
    17 = @(1317, 638, 1955)
    18 = @(95, 109, 204)
    19 = @(1440, 505, 1945)
    20 = @(24, 5, 29)
    21 = @(1661, 562, 2223)
    22 = @(224, 74, 298)
    23 = @(2279, 725, 3004)
    24 = @(10628, 2075, 12703)
    25 = @(472, 119, 591)
    27 = @(426, 218, 644)
    28 = @(342, 180, 522)
    29 = @(1831, 646, 2477)
    30 = @(12, 26, 38)
    31 = @(359, 66, 425)
    32 = @(1334, 547, 1881)
    33 = @(832, 394, 1226)
    34 = @(477, 190, 667)
    35 = @(6024, 1106, 7130)
    36 = @(683, 142, 825)
    37 = @(3196, 565, 3761)
    38 = @(147, 21, 168)
    39 = @(120, 106, 226)
    40 = @(547, 30, 577)
    41 = @(2560, 72, 2632)
    42 = @(42, 22, 64)
    43 = @(272, 10, 282)
    44 = @(395, 74, 469)
    45 = @(1330, 652, 1982)
    46 = @(2713, 973, 3686)
    47 = @(1502, 278, 1780)
    48 = @(703, 452, 1155)
    49 = @(3080, 1011, 4091)
    50 = @(538, 150, 688)
    51 = @(898, 318, 1216)
    52 = @(84, 327, 411)
    53 = @(232, 34, 266)
    54 = @(376, 535, 911)
}
foreach ($row in $rawData.Keys) {
    $vals = $rawData[$row]
    $octRaw.Cells.Item($row, 2).Value = $vals[0]
    $octRaw.Cells.Item($row, 3).Value = $vals[1]
    $octRaw.Cells.Item($row, 4).Value = $vals[2]
}

Write-Output "OctoberRaw populated; October and Yearly total will recalc."
